$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '37.230.55'
$ws.Range("E2").Value = '  +0.52%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.076.29'
$ws.Range("E3").Value = '  +0.07%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '251.70'
$ws.Range("E5").Value = '  +1.07%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.675'
$ws.Range("E6").Value = '  +4.21%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '63.50'
$ws.Range("E7").Value = '  +27.83%  '

$ws.Range("E8").Value = '  -0.07%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '61.41'
$ws.Range("E9").Value = '  +2.08%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.386'
$ws.Range("E10").Value = '  +5.46%  '

$ws.Range("E11").Value = '  +9.92%  '

$ws.Range("E12").Value = '  +2.94%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '15.96'
$ws.Range("E13").Value = '  +7.54%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.339.94'
$ws.Range("E14").Value = '  -1.31%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.829'
$ws.Range("E15").Value = '  +1.07%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.45'
$ws.Range("E16").Value = '  +8.67%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.076.39'
$ws.Range("E17").Value = '  -0.16%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '37.174.01'
$ws.Range("E18").Value = '  +0.88%  '

$ws.Range("E19").Value = '  +5.06%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0935'
$ws.Range("E20").Value = '  +14.85%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '15.26'
$ws.Range("E21").Value = '  +16.65%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.47'
$ws.Range("E22").Value = '  +7.07%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '240.41'
$ws.Range("E23").Value = '  +1.39%  '

$ws.Range("E24").Value = '  +0.11%  '

$ws.Range("E25").Value = '  -0.75%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '171.85'
$ws.Range("E26").Value = '  +2.29%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.30'
$ws.Range("E27").Value = '  +1.46%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("B28").Value = 'PancakeSwap'
$ws.Range("C28").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D28").Value = '2.06'
$ws.Range("E28").Value = '  +3.77%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("B29").Value = 'EthereumClassic'
$ws.Range("C29").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D29").Value = '20.46'
$ws.Range("E29").Value = '  -0.86%  '

$ws.Range("E30").Value = '  +4.03%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").Value = '4.72'
$ws.Range("E31").Value = '  +6.49%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").Value = '1.09'
$ws.Range("E32").Value = '  +2.74%  '

$ws.Range("E33").Value = '  +6.73%  '

$ws.Range("E34").Value = '  +10.42%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0890'
$ws.Range("E35").Value = '  -0.61%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("B36").Value = 'LidoDAOToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D36").Value = '2.33'
$ws.Range("E36").Value = '  +4.02%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("B37").Value = 'BinanceUSD'
$ws.Range("C37").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D37").Value = '1.00'
$ws.Range("E37").Value = '  -0.16%  '

$ws.Range("E39").Value = '  +26.63%  '

$ws.Range("E40").Value = '  +4.16%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '18.95'
$ws.Range("E41").Value = '  +9.44%  '

$ws.Range("E42").Value = '  +3.11%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.17'
$ws.Range("E43").Value = '  +2.81%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '98.82'
$ws.Range("E44").Value = '  +2.44%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.35'
$ws.Range("E45").Value = '  +30.16%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.80'
$ws.Range("E46").Value = '  +1.35%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.57'
$ws.Range("E47").Value = '  +17.28%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.55'
$ws.Range("E48").Value = '  +14.44%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.310.00'
$ws.Range("E49").Value = '  +1.13%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.94'

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.96'
$ws.Range("E51").Value = '  +2.68%  '
